# Update column G ("K" = strikeouts) values for giolito_lucas.xlsx save data.
# These correspond to regenerated save_data using K (strikeouts) instead of
# the previous "Strike#" (total strikes thrown) metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 4
    3 = 3
    4 = 6
    5 = 4
    6 = 8
    7 = 5
    8 = 6
    9 = 9
    10 = 8
    11 = 9
    12 = 2
    13 = 7
    14 = 3
    15 = 8
    16 = 10
    17 = 5
    18 = 1
    19 = 7
    20 = 7
    21 = 9
    22 = 9
    23 = 12
    24 = 5
    25 = 11
    26 = 7
    27 = 2
    28 = 8
    29 = 7
    30 = 0
    31 = 11
    32 = 10
    33 = 8
    34 = 9
    35 = 4
    36 = 6
    37 = 5
    38 = 2
    39 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
